# Updated cryptos list on Sun Apr 21 07:26:07 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the
# cryptocurrency rows on the active worksheet to the latest scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.004.54"
$ws.Range("E2").Value = "  +1.54%  "

$ws.Range("D3").Value = "3.179.29"
$ws.Range("E3").Value = "  +3.77%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.80"
$ws.Range("E5").Value = "  +3.15%  "

$ws.Range("E6").Value = "  +5.26%  "

$ws.Range("D8").Value = "3.178.58"
$ws.Range("E8").Value = "  +3.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  +3.24%  "

$ws.Range("E10").Value = "  +5.49%  "

$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("E12").Value = "  +2.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000274"
$ws.Range("E13").Value = "  +18.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.05"
$ws.Range("E14").Value = "  +6.32%  "

$ws.Range("D15").Value = "3.701.16"
$ws.Range("E15").Value = "  +3.82%  "

$ws.Range("D16").Value = "65.117.67"
$ws.Range("E16").Value = "  +1.62%  "

$ws.Range("D17").Value = "3.194.79"
$ws.Range("E17").Value = "  +4.22%  "

$ws.Range("E18").Value = "  +5.46%  "

$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "513.33"
$ws.Range("E20").Value = "  +7.23%  "

$ws.Range("E21").Value = "  +6.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.34"
$ws.Range("E23").Value = "  +7.04%  "

$ws.Range("E24").Value = "  +3.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.42"
$ws.Range("E25").Value = "  +3.27%  "

$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("E27").Value = "  +11.62%  "

$ws.Range("E28").Value = "  +4.63%  "

$ws.Range("E29").Value = "  +7.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.06"
$ws.Range("E30").Value = "  +6.50%  "

$ws.Range("E31").Value = "  +13.73%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("E33").Value = "  +5.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.30"
$ws.Range("E34").Value = "  +8.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.62"
$ws.Range("E35").Value = "  +6.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.72"
$ws.Range("E36").Value = "  +1.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0901"
$ws.Range("E37").Value = "  +10.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "477.38"
$ws.Range("E38").Value = "  +5.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.15"
$ws.Range("E39").Value = "  +11.40%  "

$ws.Range("E40").Value = "  +2.43%  "

$ws.Range("D42").Value = "3.062.39"
$ws.Range("E42").Value = "  +0.94%  "

$ws.Range("E43").Value = "  +2.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.288"
$ws.Range("E44").Value = "  +8.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.40"
$ws.Range("E45").Value = "  +7.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.19"
$ws.Range("E46").Value = "  +5.23%  "

$ws.Range("D47").Value = "0.0₃0615"
$ws.Range("E47").Value = "  +18.57%  "

$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("E49").Value = "  +1.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.27"
$ws.Range("E50").Value = "  +8.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.68"
$ws.Range("E51").Value = "  +1.28%  "
